$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-SaldoRow {
    param([int]$RowIndex, [string]$Conta, [string]$Nome, [double]$Saldo)

    $ws.Rows.Item($RowIndex).Insert()

    $contaCell = $ws.Cells.Item($RowIndex, 1)
    $contaCell.NumberFormat = "@"
    $contaCell.Value = $Conta
    $contaCell.ClearFormats()

    $ws.Cells.Item($RowIndex, 2).Value = $Nome
    $ws.Cells.Item($RowIndex, 3).Value = $Saldo
}

# Insert the four "moved" accounts into their new positions, bottom-to-top
# so earlier (smaller) row numbers stay valid for subsequent inserts.
Add-SaldoRow 163 "005198093" "ANA"    100
Add-SaldoRow 162 "004890544" "ASSAKO" 100
Add-SaldoRow 158 "004575632" "ADELE"  100
Add-SaldoRow 149 "002687737" "JOSE"   100

# The five original rows (with their old, negative balances) now sit four
# rows further down because of the four inserts above them (514+4 .. 518+4).
$ws.Range("A518:A522").EntireRow.Delete()
